$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round coordinates to whole numbers
$ws.Range("Q2").Value = 417636
$ws.Range("R2").Value = 6699979

# Clear the time cells (Starttid / Sluttid)
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
